$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New product data replacing the old rows 2-6 and extending through row 17.
# Columns: A=ID, B=Name, C=Price, D=Quantity, E=Category
$data = @(
    @(6,  "Czarnuch",      1,  1,    "Robol"),
    @(8,  "czarnuch ",     1,  1,    "robol"),
    @(9,  "Baton",         3,  100,  "Słodycze"),
    @(10, "Czekolada",     4,  110,  "Słodycze"),
    @(11, "Żelki Haribo",  6,  220,  "Słodycze"),
    @(12, "CocaCola",      7,  127,  "Napoje"),
    @(13, "Woda",          2,  399,  "Napoje"),
    @(14, "Fanta",         5,  75,   "Napoje"),
    @(15, "Sprite",        7,  100,  "Napoje"),
    @(16, "Pomarańcze",    2,  200,  "Owoce"),
    @(17, "Chleb biały",   4,  1000, "Pieczywo"),
    @(18, "Bułki",         1,  200,  "Pieczywo"),
    @(19, "Kalafior",      2,  22,   "Warzywa"),
    @(20, "Marchewka",     2,  36,   "Warzywa"),
    @(21, "Piwo",          2,  23,   "Alkohole"),
    @(22, "Wino",          22, 15,   "Alkohole")
)

$row = 2
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $ws.Cells.Item($row, 4).Value = $entry[3]
    $ws.Cells.Item($row, 5).Value = $entry[4]
    $row++
}
